{"js": "// Append two new paragraphs (\"Hi\" and \"ali\") after the last paragraph in the\n// document body, matching the paragraph formatting (bidi) of the paragraph\n// they follow since Word inherits pPr from the insertion point by default.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nlet hiParagraph = lastParagraph.insertParagraph(\"Hi\", Word.InsertLocation.after);\nawait context.sync();\n\nlet aliParagraph = hiParagraph.insertParagraph(\"ali\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Append two new paragraphs (\"Hi\" and \"ali\") after the last paragraph in the\n# document. InsertParagraphAfter() on the last paragraph's range creates a new\n# paragraph mark that inherits the paragraph formatting (e.g. bidi) of the\n# paragraph it follows, matching Word's normal \"press Enter at end\" behavior.\n$d = $word.ActiveDocument\n\n# Insert \"Hi\" as a new paragraph after the current last paragraph.\n$lastParagraph = $d.Paragraphs.Last\n$range = $lastParagraph.Range\n$range.InsertParagraphAfter()\n$range.Collapse(0)\n$d.Paragraphs.Last.Range.Text = \"Hi\"\n\n# Insert \"ali\" as a new paragraph after the paragraph we just added.\n$lastParagraph2 = $d.Paragraphs.Last\n$range2 = $lastParagraph2.Range\n$range2.InsertParagraphAfter()\n$range2.Collapse(0)\n$d.Paragraphs.Last.Range.Text = \"ali\"\n"}
